# Scale the "value" column (D) from the raw pushed units up to the pulled
# units by a factor of 10000 (data push/pull unit rescale), matching the
# author's "scale data push and pull" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 4).End(-4162).Row  # xlUp

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    if ($null -ne $cell.Value2) {
        # Use decimal arithmetic (rather than binary double multiplication)
        # so the rescale is an exact decimal-digit shift, matching the
        # source data re-export instead of accumulating float noise.
        # Round-trip through a string to normalize the decimal's internal
        # scale before converting back to a double (avoids a 1-ULP drift
        # that .NET's direct decimal->double conversion can introduce
        # when the decimal carries extra trailing-zero scale digits).
        $scaled = [decimal]$cell.Value2 * 10000
        $scaledNormalized = [decimal]($scaled.ToString())
        $cell.Value2 = [double]$scaledNormalized
    }
}
